$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.538.64"
$ws.Range("E2").Value = "  +5.15%  "
$ws.Range("D3").Value = "1.725.46"
$ws.Range("E3").Value = "  +4.04%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'226.22"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "'0.5384"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.2695"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "'0.06625"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("D10").Value = "'21.76"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("D11").Value = "'0.07767"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "'4.652"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "1.747.43"
$ws.Range("E13").Value = "  +5.26%  "
$ws.Range("D14").Value = "1.963.22"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").Value = "0.0₅8303"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "27.564.82"
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("D19").Value = "'225.15"
$ws.Range("E19").Value = "  +16.93%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'4.755"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").Value = "'10.75"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "'6.123"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D25").Value = "'148.11"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +11.19%  "
$ws.Range("D27").Value = "'0.1234"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "'7.429"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "'16.81"
$ws.Range("E29").Value = "  +4.96%  "
$ws.Range("D30").Value = "'0.05589"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'1.307"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").Value = "'3.590"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("D33").Value = "'3.475"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("D34").Value = "'1.670"
$ws.Range("E34").Value = "  +6.42%  "
$ws.Range("D35").Value = "'0.9663"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'2.449"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").Value = "'2.822"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").Value = "'0.5967"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").Value = "'5.905"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'0.8610"
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("D42").Value = "1.060.97"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'101.76"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "1.868.43"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("E46").Value = "  +12.90%  "
$ws.Range("D47").Value = "'59.13"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "'8.244"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").Value = "'0.4432"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "'1.006"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.05282"
$ws.Range("E51").Value = "  +0.87%  "
